$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 179, shifting rows 179:277 down to 180:278
$ws.Rows(179).Insert()

$ws.Range("A179").Value = 4
$ws.Range("B179").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C179").Value = "Los Lagos"
$ws.Range("D179").Value = 44572
$ws.Range("E179").Value = 10
$ws.Range("F179").Value = 100114013
$ws.Range("G179").Value = "Zanahoria"
$ws.Range("H179").Value = "Sin especificar"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 750
$ws.Range("K179").Value = 14000
$ws.Range("L179").Value = 14000
$ws.Range("M179").Value = 14000
$ws.Range("N179").Value = "$/saco 20 kilos"
$ws.Range("O179").Value = "Región de Ñuble"
$ws.Range("P179").Value = 700
$ws.Range("Q179").Value = 20
$ws.Range("R179").Value = "Hortaliza"
